$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2/K2: new report-type code column inserted conceptually; J2 becomes "002", K2 stays "001"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "002"
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "001"

# N2: report date updated
$ws.Range("N2").Value = "2020-06-30 00:00:00"

# Numeric financial figures updated
$ws.Range("O2").Value = 981044756.5599999
$ws.Range("P2").Value = 137813839.79
$ws.Range("Q2").Value = 121597320.37
$ws.Range("R2").Value = 104.3060627933
$ws.Range("S2").Value = 304289358.29
$ws.Range("T2").Value = -27.3850741114
$ws.Range("U2").Value = 82460862.69
$ws.Range("V2").Value = 5.9596473445
$ws.Range("W2").Value = 284700105.44
$ws.Range("X2").Value = 157234960.58
$ws.Range("Y2").Value = -12.1069287865

# Z2/AA2: values removed (now blank)
$ws.Range("Z2").ClearContents()
$ws.Range("AA2").ClearContents()

$ws.Range("AB2").Value = 696344651.12
$ws.Range("AC2").Value = 16.2189328895
$ws.Range("AD2").Value = -6.0375218332
$ws.Range("AE2").Value = -36.0102432711
$ws.Range("AF2").Value = 308.4641788622
$ws.Range("AG2").Value = 29.0200934806
